$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header in C2 was renamed from "EK" to "Buy-Price"
$ws.Range("C2").Value = "Buy-Price"

# Reflect the cell that was selected when the workbook was last saved
$ws.Range("E19").Select()
